$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1534.6364
$ws.Range("I70").Value = 1672.9
$ws.Range("K70").Value = 5018.700000000001
$ws.Range("M70").Value = -4748.700000000001
$ws.Range("H73").Value = 1534.6364
$ws.Range("I73").Value = 1672.9
$ws.Range("K73").Value = 5018.700000000001
$ws.Range("M73").Value = -4082.700000000001
$ws.Range("H80").Value = 669.5454999999999
$ws.Range("I80").Value = 478.83334
$ws.Range("J80").Value = 898.4
$ws.Range("K80").Value = 1436.50002
$ws.Range("L80").Value = 2695.2
$ws.Range("M80").Value = -438.5000199999999
$ws.Range("N80").Value = -4691.2
$ws.Range("H83").Value = 669.5454999999999
$ws.Range("I83").Value = 478.83334
$ws.Range("J83").Value = 898.4
$ws.Range("K83").Value = 4309.50006
$ws.Range("L83").Value = 8085.599999999999
$ws.Range("M83").Value = 682.4999399999997
$ws.Range("N83").Value = -18069.6
$ws.Range("H92").Value = 878.3889
$ws.Range("J92").Value = 1544
$ws.Range("L92").Value = 1544
$ws.Range("N92").Value = -4040
$ws.Range("H98").Value = 125063560
$ws.Range("I98").Value = 125063560
$ws.Range("K98").Value = 125063560
$ws.Range("M98").Value = -125062062
$ws.Range("H100").Value = 2144.0833
$ws.Range("I100").Value = 1516.125
$ws.Range("J100").Value = 3400
$ws.Range("K100").Value = 1516.125
$ws.Range("L100").Value = 3400
$ws.Range("M100").Value = -975.125
$ws.Range("N100").Value = -4482
$ws.Range("H106").Value = 8369.579
$ws.Range("I106").Value = 1918.5
$ws.Range("K106").Value = 1918.5
$ws.Range("M106").Value = -1287.5
$ws.Range("H116").Value = 6799.75
$ws.Range("I116").Value = 6078.778
$ws.Range("K116").Value = 6078.778
$ws.Range("M116").Value = -2636.778
$ws.Range("H122").Value = 125063560
$ws.Range("I122").Value = 125063560
$ws.Range("K122").Value = 375190680
$ws.Range("M122").Value = -375188230
$ws.Range("H132").Value = 2714
$ws.Range("I132").Value = 2758.25
$ws.Range("K132").Value = 8274.75
$ws.Range("M132").Value = -5744.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1499.5358
$ws.Range("I2").Value = 954.4
$ws.Range("J2").Value = 2862.375
$ws.Range("K2").Value = 954.4
$ws.Range("L2").Value = 2862.375
$ws.Range("M2").Value = -841.4
$ws.Range("N2").Value = -3088.375
$ws.Range("H32").Value = 16670238
$ws.Range("I32").Value = 18520856
$ws.Range("K32").Value = 18520856
$ws.Range("M32").Value = -18520569
$ws.Range("H45").Value = 2383.2942
$ws.Range("I45").Value = 1750.25
$ws.Range("K45").Value = 1750.25
$ws.Range("M45").Value = -1373.25
$ws.Range("H61").Value = 13192523
$ws.Range("I61").Value = 15158460
$ws.Range("K61").Value = 15158460
$ws.Range("M61").Value = -15158248
$ws.Range("H63").Value = 6261.6875
$ws.Range("I63").Value = 4743.7
$ws.Range("K63").Value = 4743.7
$ws.Range("M63").Value = -4057.7
$ws.Range("H66").Value = 6261.6875
$ws.Range("I66").Value = 4743.7
$ws.Range("K66").Value = 23718.5
$ws.Range("M66").Value = -20286.5
$ws.Range("H74").Value = 9624649
$ws.Range("I74").Value = 19232668
$ws.Range("J74").Value = 16629.309
$ws.Range("K74").Value = 19232668
$ws.Range("L74").Value = 16629.309
$ws.Range("M74").Value = -19231794
$ws.Range("N74").Value = -18377.309
$ws.Range("H77").Value = 9624649
$ws.Range("I77").Value = 19232668
$ws.Range("J77").Value = 16629.309
$ws.Range("K77").Value = 96163340
$ws.Range("L77").Value = 83146.54500000001
$ws.Range("M77").Value = -96158972
$ws.Range("N77").Value = -91882.54500000001
$ws.Range("H102").Value = 32620.666
$ws.Range("I102").Value = 36073.25
$ws.Range("K102").Value = 36073.25
$ws.Range("M102").Value = -34451.25
$ws.Range("H106").Value = 39684.832
$ws.Range("J106").Value = 39684.832
$ws.Range("L106").Value = 39684.832
$ws.Range("N106").Value = -42208.832
$ws.Range("H116").Value = 1499.5358
$ws.Range("I116").Value = 954.4
$ws.Range("J116").Value = 2862.375
$ws.Range("K116").Value = 954.4
$ws.Range("L116").Value = 2862.375
$ws.Range("M116").Value = 1339.6
$ws.Range("N116").Value = -7450.375
$ws.Range("H122").Value = 2706.75
$ws.Range("I122").Value = 1068.7142
$ws.Range("K122").Value = 3206.1426
$ws.Range("M122").Value = -756.1425999999997
$ws.Range("H136").Value = 13192523
$ws.Range("I136").Value = 15158460
$ws.Range("K136").Value = 45475380
$ws.Range("M136").Value = -45472830

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1499.5358
$ws.Range("I3").Value = 954.4
$ws.Range("J3").Value = 2862.375
$ws.Range("K3").Value = 954.4
$ws.Range("L3").Value = 2862.375
$ws.Range("M3").Value = -840.4
$ws.Range("N3").Value = -3090.375
$ws.Range("H94").Value = 978.96295
$ws.Range("I94").Value = 1114.8889
$ws.Range("K94").Value = 1114.8889
$ws.Range("M94").Value = -663.8888999999999
$ws.Range("H98").Value = 57013.668
$ws.Range("J98").Value = 57013.668
$ws.Range("L98").Value = 57013.668
$ws.Range("N98").Value = -63003.668
$ws.Range("H134").Value = 82405.766
$ws.Range("I134").Value = 1124.7
$ws.Range("K134").Value = 3374.1
$ws.Range("M134").Value = -839.1000000000004

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3117.5
$ws.Range("I16").Value = 1141
$ws.Range("K16").Value = 1141
$ws.Range("M16").Value = -854
$ws.Range("H31").Value = 381573.06
$ws.Range("I31").Value = 4726.4
$ws.Range("J31").Value = 1066748.8
$ws.Range("K31").Value = 4726.4
$ws.Range("L31").Value = 1066748.8
$ws.Range("M31").Value = -4431.4
$ws.Range("N31").Value = -1067338.8
$ws.Range("H34").Value = 381573.06
$ws.Range("I34").Value = 4726.4
$ws.Range("J34").Value = 1066748.8
$ws.Range("K34").Value = 4726.4
$ws.Range("L34").Value = 1066748.8
$ws.Range("M34").Value = -4524.4
$ws.Range("N34").Value = -1067152.8
$ws.Range("H107").Value = 1912.6875
$ws.Range("J107").Value = 2315.5
$ws.Range("L107").Value = 2315.5
$ws.Range("N107").Value = -6155.5
$ws.Range("H113").Value = 3117.5
$ws.Range("I113").Value = 1141
$ws.Range("K113").Value = 1141
$ws.Range("M113").Value = 1029

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 40434.812
$ws.Range("I44").Value = 40434.812
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 121304.436
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -120906.436
$ws.Range("H55").Value = 8957.286
$ws.Range("J55").Value = 7995
$ws.Range("L55").Value = 23985
$ws.Range("N55").Value = -24339
$ws.Range("H58").Value = 1824.5
$ws.Range("J58").Value = 1499.3334
$ws.Range("L58").Value = 4498.0002
$ws.Range("N58").Value = -4754.0002
$ws.Range("H122").Value = 948.4545000000001
$ws.Range("J122").Value = 948.4545000000001
$ws.Range("L122").Value = 8536.0905
$ws.Range("N122").Value = -13436.0905
$ws.Range("N44").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3872.889
$ws.Range("J113").Value = 4005.625
$ws.Range("L113").Value = 4005.625
$ws.Range("N113").Value = -8345.625
$ws.Range("H123").Value = 29224
$ws.Range("J123").Value = 29224
$ws.Range("L123").Value = 29224
$ws.Range("N123").Value = -34124

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14360866
$ws.Range("I7").Value = 22224566
$ws.Range("K7").Value = 22224566
$ws.Range("M7").Value = -22224454
$ws.Range("H82").Value = 1071.5
$ws.Range("I82").Value = 606.75
$ws.Range("K82").Value = 606.75
$ws.Range("M82").Value = -245.75
$ws.Range("H85").Value = 1071.5
$ws.Range("I85").Value = 606.75
$ws.Range("K85").Value = 606.75
$ws.Range("M85").Value = 641.25
$ws.Range("H126").Value = 14360866
$ws.Range("I126").Value = 22224566
$ws.Range("K126").Value = 66673698
$ws.Range("M126").Value = -66671228
$ws.Range("H127").Value = 155950
$ws.Range("J127").Value = 155950
$ws.Range("L127").Value = 155950
$ws.Range("N127").Value = -165870

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6673980
$ws.Range("I62").Value = 7733.3335
$ws.Range("J62").Value = 16673350
$ws.Range("K62").Value = 7733.3335
$ws.Range("L62").Value = 16673350
$ws.Range("M62").Value = -7109.3335
$ws.Range("N62").Value = -16674598
$ws.Range("H65").Value = 6673980
$ws.Range("I65").Value = 7733.3335
$ws.Range("J65").Value = 16673350
$ws.Range("K65").Value = 38666.6675
$ws.Range("L65").Value = 83366750
$ws.Range("M65").Value = -35546.6675
$ws.Range("N65").Value = -83372990
$ws.Range("H109").Value = 55800
$ws.Range("J109").Value = 55800
$ws.Range("L109").Value = 55800
$ws.Range("N109").Value = -58574
$ws.Range("H122").Value = 9756.143
$ws.Range("I122").Value = 3582.25
$ws.Range("J122").Value = 17988
$ws.Range("K122").Value = 10746.75
$ws.Range("L122").Value = 53964
$ws.Range("M122").Value = -8296.75
$ws.Range("N122").Value = -58864
$ws.Range("H132").Value = 1764.359
$ws.Range("I132").Value = 1832.5358
$ws.Range("K132").Value = 5497.607400000001
$ws.Range("M132").Value = -2967.607400000001
